$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Cells.Item(3, 4).Value = 44511
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 7000
$ws.Cells.Item(3, 13).Value = 7500
$ws.Cells.Item(3, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(3, 16).Value = 300

# Row 4
$ws.Cells.Item(4, 4).Value = 44509
$ws.Cells.Item(4, 11).Value = 8000
$ws.Cells.Item(4, 12).Value = 9000
$ws.Cells.Item(4, 13).Value = 8500
$ws.Cells.Item(4, 15).Value = 'Región del Maule'
$ws.Cells.Item(4, 16).Value = 340

# Row 5
$ws.Cells.Item(5, 4).Value = 44495
$ws.Cells.Item(5, 11).Value = 8000
$ws.Cells.Item(5, 12).Value = 9000
$ws.Cells.Item(5, 13).Value = 8500
$ws.Cells.Item(5, 15).Value = 'Región del Maule'
$ws.Cells.Item(5, 16).Value = 340

# Row 6
$ws.Cells.Item(6, 4).Value = 44162
$ws.Cells.Item(6, 10).Value = 80
$ws.Cells.Item(6, 11).Value = 7000
$ws.Cells.Item(6, 12).Value = 8000
$ws.Cells.Item(6, 13).Value = 7562
$ws.Cells.Item(6, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(6, 16).Value = 302

# Row 8
$ws.Cells.Item(8, 4).Value = 44161
$ws.Cells.Item(8, 10).Value = 53
$ws.Cells.Item(8, 11).Value = 6500
$ws.Cells.Item(8, 12).Value = 7000
$ws.Cells.Item(8, 13).Value = 6764
$ws.Cells.Item(8, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(8, 16).Value = 271

# Row 9
$ws.Cells.Item(9, 4).Value = 44166
$ws.Cells.Item(9, 10).Value = 56
$ws.Cells.Item(9, 11).Value = 7500
$ws.Cells.Item(9, 12).Value = 8000
$ws.Cells.Item(9, 13).Value = 7804
$ws.Cells.Item(9, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(9, 16).Value = 312

# Row 10
$ws.Cells.Item(10, 4).Value = 44530

# Row 11
$ws.Cells.Item(11, 4).Value = 44491
$ws.Cells.Item(11, 10).Value = 60

# Row 12
$ws.Cells.Item(12, 4).Value = 44167
$ws.Cells.Item(12, 10).Value = 60
$ws.Cells.Item(12, 11).Value = 8000
$ws.Cells.Item(12, 12).Value = 9000
$ws.Cells.Item(12, 13).Value = 8500
$ws.Cells.Item(12, 15).Value = 'Región del Maule'
$ws.Cells.Item(12, 16).Value = 340

# Row 13
$ws.Cells.Item(13, 4).Value = 44553
$ws.Cells.Item(13, 10).Value = 100
$ws.Cells.Item(13, 11).Value = 6500
$ws.Cells.Item(13, 13).Value = 6750
$ws.Cells.Item(13, 16).Value = 270

# Row 14
$ws.Cells.Item(14, 4).Value = 44526
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 6000
$ws.Cells.Item(14, 13).Value = 6500
$ws.Cells.Item(14, 16).Value = 260

# Row 15
$ws.Cells.Item(15, 4).Value = 44487
$ws.Cells.Item(15, 10).Value = 30
$ws.Cells.Item(15, 11).Value = 8000
$ws.Cells.Item(15, 12).Value = 8000
$ws.Cells.Item(15, 13).Value = 8000
$ws.Cells.Item(15, 15).Value = 'Región del Maule'
$ws.Cells.Item(15, 16).Value = 320

# Row 16
$ws.Cells.Item(16, 4).Value = 44487
$ws.Cells.Item(16, 9).Value = 'Segunda'
$ws.Cells.Item(16, 10).Value = 30
$ws.Cells.Item(16, 11).Value = 9000
$ws.Cells.Item(16, 12).Value = 9000
$ws.Cells.Item(16, 13).Value = 9000
$ws.Cells.Item(16, 15).Value = 'Región del Maule'
$ws.Cells.Item(16, 16).Value = 360

# Row 17
$ws.Cells.Item(17, 4).Value = 44537
$ws.Cells.Item(17, 10).Value = 60
$ws.Cells.Item(17, 11).Value = 6500
$ws.Cells.Item(17, 12).Value = 7000
$ws.Cells.Item(17, 13).Value = 6750
$ws.Cells.Item(17, 16).Value = 270

# Row 18
$ws.Cells.Item(18, 4).Value = 44482
$ws.Cells.Item(18, 10).Value = 120

# Row 19
$ws.Cells.Item(19, 4).Value = 44160
$ws.Cells.Item(19, 10).Value = 80
$ws.Cells.Item(19, 11).Value = 6500
$ws.Cells.Item(19, 13).Value = 6688
$ws.Cells.Item(19, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(19, 16).Value = 268

# Row 20
$ws.Cells.Item(20, 4).Value = 44473
$ws.Cells.Item(20, 10).Value = 60
$ws.Cells.Item(20, 11).Value = 9500
$ws.Cells.Item(20, 12).Value = 10000
$ws.Cells.Item(20, 13).Value = 9750
$ws.Cells.Item(20, 15).Value = 'Región del Maule'
$ws.Cells.Item(20, 16).Value = 390

# Row 21
$ws.Cells.Item(21, 4).Value = 44536
$ws.Cells.Item(21, 10).Value = 80
$ws.Cells.Item(21, 13).Value = 6750
$ws.Cells.Item(21, 15).Value = 'Provincia de Diguillín'

# Row 22
$ws.Cells.Item(22, 4).Value = 44529
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 6000
$ws.Cells.Item(22, 12).Value = 7000
$ws.Cells.Item(22, 13).Value = 6500
$ws.Cells.Item(22, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(22, 16).Value = 260

# Row 23
$ws.Cells.Item(23, 4).Value = 44517
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 11).Value = 6000
$ws.Cells.Item(23, 12).Value = 7000
$ws.Cells.Item(23, 13).Value = 6500
$ws.Cells.Item(23, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(23, 16).Value = 260

# Row 24
$ws.Cells.Item(24, 4).Value = 44518
$ws.Cells.Item(24, 10).Value = 60
$ws.Cells.Item(24, 11).Value = 6000
$ws.Cells.Item(24, 13).Value = 6500
$ws.Cells.Item(24, 16).Value = 260

# Row 25
$ws.Cells.Item(25, 4).Value = 44540
$ws.Cells.Item(25, 10).Value = 100
$ws.Cells.Item(25, 11).Value = 6500
$ws.Cells.Item(25, 13).Value = 6750
$ws.Cells.Item(25, 16).Value = 270

# Row 26
$ws.Cells.Item(26, 4).Value = 44488
$ws.Cells.Item(26, 10).Value = 60

# Row 27
$ws.Cells.Item(27, 4).Value = 44466
$ws.Cells.Item(27, 10).Value = 60
$ws.Cells.Item(27, 11).Value = 11000
$ws.Cells.Item(27, 12).Value = 12000
$ws.Cells.Item(27, 13).Value = 11500
$ws.Cells.Item(27, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(27, 16).Value = 460

# Row 28
$ws.Cells.Item(28, 4).Value = 44524
$ws.Cells.Item(28, 10).Value = 100
$ws.Cells.Item(28, 11).Value = 6000
$ws.Cells.Item(28, 13).Value = 6500
$ws.Cells.Item(28, 16).Value = 260

# Row 29
$ws.Cells.Item(29, 4).Value = 44512
$ws.Cells.Item(29, 10).Value = 100
$ws.Cells.Item(29, 11).Value = 7000
$ws.Cells.Item(29, 12).Value = 8000
$ws.Cells.Item(29, 13).Value = 7500
$ws.Cells.Item(29, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(29, 16).Value = 300

# Row 30
$ws.Cells.Item(30, 4).Value = 44519
$ws.Cells.Item(30, 10).Value = 80
$ws.Cells.Item(30, 11).Value = 6000
$ws.Cells.Item(30, 12).Value = 7000
$ws.Cells.Item(30, 13).Value = 6500
$ws.Cells.Item(30, 16).Value = 260

# Row 31
$ws.Cells.Item(31, 4).Value = 44523
$ws.Cells.Item(31, 10).Value = 80
$ws.Cells.Item(31, 11).Value = 6000
$ws.Cells.Item(31, 12).Value = 7000
$ws.Cells.Item(31, 13).Value = 6500
$ws.Cells.Item(31, 16).Value = 260

# Row 32
$ws.Cells.Item(32, 4).Value = 44566
$ws.Cells.Item(32, 10).Value = 60
$ws.Cells.Item(32, 12).Value = 7500
$ws.Cells.Item(32, 13).Value = 7250
$ws.Cells.Item(32, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(32, 16).Value = 290

# Row 33
$ws.Cells.Item(33, 4).Value = 44476
$ws.Cells.Item(33, 10).Value = 160
$ws.Cells.Item(33, 11).Value = 7500
$ws.Cells.Item(33, 12).Value = 8000
$ws.Cells.Item(33, 13).Value = 7750
$ws.Cells.Item(33, 15).Value = 'Región del Maule'
$ws.Cells.Item(33, 16).Value = 310

# Row 34
$ws.Cells.Item(34, 4).Value = 44533
$ws.Cells.Item(34, 10).Value = 80
$ws.Cells.Item(34, 11).Value = 6500
$ws.Cells.Item(34, 12).Value = 7000
$ws.Cells.Item(34, 13).Value = 6750
$ws.Cells.Item(34, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(34, 16).Value = 270

# Row 35
$ws.Cells.Item(35, 4).Value = 44489
$ws.Cells.Item(35, 9).Value = 'Primera'
$ws.Cells.Item(35, 10).Value = 60
$ws.Cells.Item(35, 11).Value = 8000
$ws.Cells.Item(35, 13).Value = 8500
$ws.Cells.Item(35, 16).Value = 340

# Row 36
$ws.Cells.Item(36, 4).Value = 44159
$ws.Cells.Item(36, 10).Value = 42
$ws.Cells.Item(36, 11).Value = 6500
$ws.Cells.Item(36, 12).Value = 7000
$ws.Cells.Item(36, 13).Value = 6738
$ws.Cells.Item(36, 16).Value = 270

# Row 37
$ws.Cells.Item(37, 4).Value = 44515
$ws.Cells.Item(37, 10).Value = 100
$ws.Cells.Item(37, 11).Value = 7000
$ws.Cells.Item(37, 12).Value = 8000
$ws.Cells.Item(37, 13).Value = 7500
$ws.Cells.Item(37, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(37, 16).Value = 300

# Row 38
$ws.Cells.Item(38, 4).Value = 44504
$ws.Cells.Item(38, 10).Value = 60
$ws.Cells.Item(38, 11).Value = 8000
$ws.Cells.Item(38, 12).Value = 9000
$ws.Cells.Item(38, 13).Value = 8500
$ws.Cells.Item(38, 15).Value = 'Región del Maule'
$ws.Cells.Item(38, 16).Value = 340

# Row 39
$ws.Cells.Item(39, 4).Value = 44484
$ws.Cells.Item(39, 10).Value = 30
$ws.Cells.Item(39, 11).Value = 8500
$ws.Cells.Item(39, 12).Value = 9000
$ws.Cells.Item(39, 13).Value = 8750
$ws.Cells.Item(39, 15).Value = 'Región del Maule'
$ws.Cells.Item(39, 16).Value = 350

# Row 40
$ws.Cells.Item(40, 4).Value = 44516
$ws.Cells.Item(40, 10).Value = 100
$ws.Cells.Item(40, 11).Value = 7000
$ws.Cells.Item(40, 12).Value = 8000
$ws.Cells.Item(40, 13).Value = 7500
$ws.Cells.Item(40, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(40, 16).Value = 300
